$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N ("Late"), shifting
# Late/Heading/Outstanding columns one position to the right (N->O, O->P, P->Q).
$ws.Columns("N").Insert() | Out-Null

# Give the freshly inserted column a manual width (it is not an auto
# "bestFit" column since it was added by hand, unlike its neighbours).
$ws.Columns("N").ColumnWidth = 10.25

# Activate the "Repayment Schedule" sheet (it becomes the active tab,
# replacing "Transactions") and select cell O8 on it.
$ws.Activate() | Out-Null
$ws.Range("O8").Select() | Out-Null
